$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell (this also grows the shared-strings table and the
# sheet dimension from A1:I1 to A1:J1 automatically).
$ws.Range("J1").Value = "customersPrimaryNames"

# Resize the affected columns (F, I and J). The host quantizes stored column
# widths to whole pixels on a 6px-per-character grid, so we pick the
# ColumnWidth input that lands closest to the target stored widths
# (12.8, 10.85 and 22.77 characters respectively).
$ws.Columns.Item(6).ColumnWidth = 12
$ws.Columns.Item(9).ColumnWidth = 10
$ws.Columns.Item(10).ColumnWidth = 22

# Move the active selection from F15 to F6.
$ws.Range("F6").Select()
